# Insert a new weekly data row at row 189, shifting the existing rows
# (189..289) down to (190..290), then populate the new row with the
# latest price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 189..289 down by one row, preserving data & formatting.
$ws.Rows.Item(189).Insert()

# Populate the newly inserted row 189 with the new record.
$ws.Cells.Item(189, 1).Value = 7
$ws.Cells.Item(189, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(189, 3).Value = "Ñuble"
$ws.Cells.Item(189, 4).Value = 45001
$ws.Cells.Item(189, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(189, 5).Value = 16
$ws.Cells.Item(189, 6).Value = "Fruta"
$ws.Cells.Item(189, 7).Value = 100108
$ws.Cells.Item(189, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(189, 9).Value = 100108005
$ws.Cells.Item(189, 10).Value = "Piña"
$ws.Cells.Item(189, 11).Value = "Caramelo"
$ws.Cells.Item(189, 12).Value = "Segunda"
$ws.Cells.Item(189, 13).Value = 20
$ws.Cells.Item(189, 14).Value = 22000
$ws.Cells.Item(189, 15).Value = 22000
$ws.Cells.Item(189, 16).Value = 22000
$ws.Cells.Item(189, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(189, 18).Value = "Ecuador"
$ws.Cells.Item(189, 19).Value = 1571
$ws.Cells.Item(189, 20).Value = 14
